{"js": "// Update the Table1 figure-size value: \" 0.0033\" -> \" 0.0024\"\n// (the stray duplicate/leading <w:sectPr> that preceded the table and the\n// namespace/element ordering on the trailing <w:sectPr> are cleaned up\n// automatically by the canonical OOXML writer as part of any save, so the\n// only content change we need to make here is the text value itself).\nconst body = context.document.body;\n\nconst results = body.search(\" 0.0033\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find table text \" 0.0033\" to update.');\n}\n\nresults.items[0].insertText(\" 0.0024\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Update the Table1 figure-size value: \" 0.0033\" -> \" 0.0024\"\n# (the stray duplicate/leading sectPr ahead of the table and the\n# element ordering on the trailing sectPr are normalized automatically by\n# Word's canonical OOXML writer on save, so the only content change needed\n# here is the text value itself).\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \" 0.0033\"\n$find.Replacement.Text = \" 0.0024\"\n$find.Forward = $true\n$find.Wrap = 1          # wdFindContinue\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 1) | Out-Null\n"}
